# Weekly refresh: insert two new price rows (current week) right after the
# existing row 123, pushing the rest of the historical rows down by two.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 124:125 - everything from the old row 124 onward
# (previously ending at row 149) shifts down to 126:151.
$ws.Rows("124:125").Insert()

# New row 124 - "Primera" quality entry for the latest reporting date.
$ws.Range("A124").Value = 7
$ws.Range("B124").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C124").Value = "Ñuble"
$ws.Range("D124").Value = 44504
$ws.Range("E124").Value = 16
$ws.Range("F124").Value = "Fruta"
$ws.Range("G124").Value = 100101
$ws.Range("H124").Value = "Berries"
$ws.Range("I124").Value = 100112025
$ws.Range("J124").Value = "Frutilla"
$ws.Range("K124").Value = "Sin especificar"
$ws.Range("L124").Value = "Primera"
$ws.Range("M124").Value = 160
$ws.Range("N124").Value = 6000
$ws.Range("O124").Value = 6500
$ws.Range("P124").Value = 6250
$ws.Range("Q124").Value = "$/caja 7 kilos"
$ws.Range("R124").Value = "Provincia de Diguillín"
$ws.Range("S124").Value = 893
$ws.Range("T124").Value = 7

# New row 125 - "Segunda" quality entry for the same latest reporting date.
$ws.Range("A125").Value = 7
$ws.Range("B125").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C125").Value = "Ñuble"
$ws.Range("D125").Value = 44504
$ws.Range("E125").Value = 16
$ws.Range("F125").Value = "Fruta"
$ws.Range("G125").Value = 100101
$ws.Range("H125").Value = "Berries"
$ws.Range("I125").Value = 100112025
$ws.Range("J125").Value = "Frutilla"
$ws.Range("K125").Value = "Sin especificar"
$ws.Range("L125").Value = "Segunda"
$ws.Range("M125").Value = 120
$ws.Range("N125").Value = 5000
$ws.Range("O125").Value = 5500
$ws.Range("P125").Value = 5250
$ws.Range("Q125").Value = "$/caja 7 kilos"
$ws.Range("R125").Value = "Provincia de Diguillín"
$ws.Range("S125").Value = 750
$ws.Range("T125").Value = 7
